$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 22.04848781383541
$ws.Range("C2").Value = 7.086576677024649
$ws.Range("D2").Value = 9.7727564605112
$ws.Range("E2").Value = 10.02013192933921
$ws.Range("F2").Value = 51.48609614823078
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("L2").Value = 10.38116335073147
$ws.Range("B3").Value = 21.76138917492582
$ws.Range("C3").Value = 6.6278303474271
$ws.Range("D3").Value = 9.654196310774976
$ws.Range("E3").Value = 10.0179263290129
$ws.Range("F3").Value = 50.3393507490849
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("L3").Value = 10.37744982034068
$ws.Range("B4").Value = 21.59331492376993
$ws.Range("C4").Value = 6.329235690193406
$ws.Range("D4").Value = 9.581002176144137
$ws.Range("E4").Value = 10.01709169825472
$ws.Range("F4").Value = 49.63020320537576
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("L4").Value = 10.37764010772768
$ws.Range("B5").Value = 21.52698135237781
$ws.Range("C5").Value = 6.203291610621427
$ws.Range("D5").Value = 9.55108964639291
$ws.Range("E5").Value = 10.0168817615722
$ws.Range("F5").Value = 49.34028504998999
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("L5").Value = 10.37833839773075
$ws.Range("B6").Value = 21.51609971108215
$ws.Range("C6").Value = 6.182120904925941
$ws.Range("D6").Value = 9.54611796555937
$ws.Range("E6").Value = 10.01685475019237
$ws.Range("F6").Value = 49.29209781294151
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("L6").Value = 10.37849181934008
$ws.Range("B7").Value = 21.59241146710758
$ws.Range("C7").Value = 6.327554441327651
$ws.Range("D7").Value = 9.58059909153673
$ws.Range("E7").Value = 10.01708834051207
$ws.Range("F7").Value = 49.62629660586732
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("L7").Value = 10.37764701258595
$ws.Range("B8").Value = 21.9478516463581
$ws.Range("C8").Value = 6.931907290069591
$ws.Range("D8").Value = 9.731968318559584
$ws.Range("E8").Value = 10.01926330444802
$ws.Range("F8").Value = 51.09194122699876
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("L8").Value = 10.37937001928296
$ws.Range("B9").Value = 22.70556283401938
$ws.Range("C9").Value = 7.983254905178049
$ws.Range("D9").Value = 10.02500069203513
$ws.Range("E9").Value = 10.0276734980628
$ws.Range("F9").Value = 53.91201034155094
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("L9").Value = 10.4023569889279
$ws.Range("B10").Value = 23.29308280119709
$ws.Range("C10").Value = 8.674878739815446
$ws.Range("D10").Value = 10.23711248379187
$ws.Range("E10").Value = 10.03641090551793
$ws.Range("F10").Value = 55.93310628958029
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("L10").Value = 10.4311880355025
$ws.Range("B11").Value = 23.56570214304285
$ws.Range("C11").Value = 8.972149898703801
$ws.Range("D11").Value = 10.33273607416815
$ws.Range("E11").Value = 10.04094794405023
$ws.Range("F11").Value = 56.83805544552235
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("L11").Value = 10.44688716086001
$ws.Range("B12").Value = 23.66959241502695
$ws.Range("C12").Value = 9.082238192239862
$ws.Range("D12").Value = 10.36880640503173
$ws.Range("E12").Value = 10.04274735343877
$ws.Range("F12").Value = 57.17839089387275
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("L12").Value = 10.45320200225986
$ws.Range("B13").Value = 23.6471904552624
$ws.Range("C13").Value = 9.058638872222669
$ws.Range("D13").Value = 10.36104449596189
$ws.Range("E13").Value = 10.04235619390333
$ws.Range("F13").Value = 57.10520215084104
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("L13").Value = 10.45182556552832
$ws.Range("B14").Value = 23.57423677825718
$ws.Range("C14").Value = 8.981256590749267
$ws.Range("D14").Value = 10.33570650078831
$ws.Range("E14").Value = 10.04109434990611
$ws.Range("F14").Value = 56.86610361032422
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("L14").Value = 10.44739928346166
$ws.Range("B15").Value = 23.52963249567538
$ws.Range("C15").Value = 8.933534903139748
$ws.Range("D15").Value = 10.32016752619046
$ws.Range("E15").Value = 10.04033203942845
$ws.Range("F15").Value = 56.71933521663816
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("L15").Value = 10.44473617834687
$ws.Range("B16").Value = 23.27536439424669
$ws.Range("C16").Value = 8.655103335785927
$ws.Range("D16").Value = 10.23084447217589
$ws.Range("E16").Value = 10.03612575538853
$ws.Range("F16").Value = 55.87365129645926
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("L16").Value = 10.43021392287692
$ws.Range("B17").Value = 23.12066850303994
$ws.Range("C17").Value = 8.479859689487453
$ws.Range("D17").Value = 10.17581442349048
$ws.Range("E17").Value = 10.03368970264779
$ws.Range("F17").Value = 55.35095063802557
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("L17").Value = 10.42196571995115
$ws.Range("B18").Value = 23.03220224888796
$ws.Range("C18").Value = 8.377430098039792
$ws.Range("D18").Value = 10.14408189171638
$ws.Range("E18").Value = 10.03234145325281
$ws.Range("F18").Value = 55.0489612424598
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("L18").Value = 10.41746488674031
$ws.Range("B19").Value = 23.00234032116643
$ws.Range("C19").Value = 8.34246814371603
$ws.Range("D19").Value = 10.13332441950556
$ws.Range("E19").Value = 10.03189403105694
$ws.Range("F19").Value = 54.94649021206109
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("L19").Value = 10.41598281440406
$ws.Range("B20").Value = 23.13708411956545
$ws.Range("C20").Value = 8.498683700539198
$ws.Range("D20").Value = 10.18168093219139
$ws.Range("E20").Value = 10.03394354587316
$ws.Range("F20").Value = 55.40673431633117
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("L20").Value = 10.42281858354289
$ws.Range("B21").Value = 23.59564816571286
$ws.Range("C21").Value = 9.004052890659551
$ws.Range("D21").Value = 10.34315281294248
$ws.Range("E21").Value = 10.04146277310705
$ws.Range("F21").Value = 56.93639846475115
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("L21").Value = 10.44868936495747
$ws.Range("B22").Value = 23.8991134456424
$ws.Range("C22").Value = 9.319881258566401
$ws.Range("D22").Value = 10.44785956249913
$ws.Range("E22").Value = 10.04685129296865
$ws.Range("F22").Value = 57.92231747493432
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("L22").Value = 10.46775259402398
$ws.Range("B23").Value = 23.73684111789489
$ws.Range("C23").Value = 9.152636257175583
$ws.Range("D23").Value = 10.39205609605437
$ws.Range("E23").Value = 10.04393179150411
$ws.Range("F23").Value = 57.39746157877453
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("L23").Value = 10.4573816232531
$ws.Range("B24").Value = 23.12966114520576
$ws.Range("C24").Value = 8.490178597992188
$ws.Range("D24").Value = 10.17902897840229
$ws.Range("E24").Value = 10.03382862057738
$ws.Range("F24").Value = 55.38151911022086
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("L24").Value = 10.42243225260898
$ws.Range("B25").Value = 22.49471132685935
$ws.Range("C25").Value = 7.713085235959534
$ws.Range("D25").Value = 9.946223733471873
$ws.Range("E25").Value = 10.02495125711189
$ws.Range("F25").Value = 53.15679566544232
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("L25").Value = 10.39403885515954
